$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "WTI"

# B3 holds a plain text date-like string ("2019-01-07"). A bare .Value
# assignment lets Excel's type inference turn it into a real date serial,
# so force Text format first, assign the literal string, then clear the
# formatting override again so the cell's style matches the original
# (unstyled) cell instead of picking up a new "@" number format.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2019-01-07"
$ws.Range("B3").ClearFormats()

$ws.Range("B4").Value = 48.27
